$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J (year 2020) ---------------------------------------------

# J3: empty cell, just picks up the thick bottom border used across row 3.
$ws.Range("J3").Borders(9).LineStyle = 1
$ws.Range("J3").Borders(9).Weight = -4138

# J4: header year value, formatted like I4.
$ws.Range("I4").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("J4").Value = 2020

# J5
$ws.Range("I5").Copy()
$ws.Range("J5").PasteSpecial(-4122)
$ws.Range("J5").Value = 370

# J6
$ws.Range("I6").Copy()
$ws.Range("J6").PasteSpecial(-4122)
$ws.Range("J6").Value = 5

# J7
$ws.Range("I7").Copy()
$ws.Range("J7").PasteSpecial(-4122)
$ws.Range("J7").Value = 5

# J8
$ws.Range("I8").Copy()
$ws.Range("J8").PasteSpecial(-4122)
$ws.Range("J8").Value = 20

# J9
$ws.Range("I9").Copy()
$ws.Range("J9").PasteSpecial(-4122)
$ws.Range("J9").Value = 19

# J10
$ws.Range("I10").Copy()
$ws.Range("J10").PasteSpecial(-4122)
$ws.Range("J10").Value = 73

# --- Column I corrections ---------------------------------------------------

$ws.Range("I8").Value = 42
$ws.Range("I9").Value = 30
$ws.Range("I10").Value = 62

Write-Output "done"
